# Updates to pass ConvertLeadToAccount test case.
# Adds a new "Opp" worksheet (between NewAccounts and VisualforceToLWC)
# containing a small Opportunity data table used by the test, and makes
# it the active/selected sheet.

$wb = $excel.ActiveWorkbook

# Insert the new sheet right after "NewAccounts" so the sheet order becomes
# NewAccounts, Opp, VisualforceToLWC (matches target sheetId/r:id layout).
$afterSheet = $wb.Worksheets.Item("NewAccounts")
$ws = $wb.Worksheets.Add($null, $afterSheet)
$ws.Name = "Opp"

# Header row
$ws.Range("A1").Value = "OppName"
$ws.Range("B1").Value = "Probability"
$ws.Range("C1").Value = "Stage"
$ws.Range("D1").Value = "CloseDate"

# Data row
$ws.Range("A2").Value = "TestAutomation"
$ws.Range("B2").Value = 20
$ws.Range("C2").Value = "Needs Analysis"

# Column sizing to roughly match the authored bestFit widths.
$ws.Columns.Item(1).ColumnWidth = 13.5
$ws.Columns.Item(2).ColumnWidth = 12.666666666666666
$ws.Columns.Item(3).ColumnWidth = 12.666666666666666
$ws.Columns.Item(4).ColumnWidth = 8.5

# Leave the selection/active cell on D2, and make this the active sheet/tab.
[void]$ws.Range("D2").Select()
[void]$ws.Activate()
